# Updates cryptos list values (price / volume%) for Wed Sep 27 18:28:27 UTC 2023 run.
# Cells are plain text (inlineStr) in the source sheet, so assign literal
# strings via Range.Value to avoid Excel's automatic number/percentage
# coercion (e.g. "1.00" staying "1.00", not becoming 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '26.214.37'
    'E2' = '  +0.10%  '
    'D3' = '1.591.08'
    'E3' = '  +0.32%  '
    'E4' = '  -0.04%  '
    'D5' = '211.75'
    'E5' = '  -0.18%  '
    'E6' = '  -0.47%  '
    'E7' = '  -0.01%  '
    'E8' = '  -0.27%  '
    'E9' = '  +0.22%  '
    'D10' = '18.86'
    'E10' = '  -1.73%  '
    'D12' = '1.813.63'
    'E12' = '  +0.25%  '
    'D13' = '1.609.86'
    'E13' = '  +1.33%  '
    'E14' = '  -0.13%  '
    'E15' = '  -2.59%  '
    'E16' = '  -0.76%  '
    'D17' = '26.189.55'
    'E17' = '  +0.01%  '
    'D18' = '228.20'
    'E18' = '  +6.52%  '
    'D19' = '0.0₃0719'
    'E19' = '  -0.73%  '
    'D20' = '7.55'
    'E20' = '  +3.85%  '
    'E21' = '  +0.03%  '
    'D22' = '4.23'
    'E22' = '  -0.31%  '
    'D23' = '2.15'
    'E23' = '  +1.45%  '
    'D24' = '8.88'
    'E24' = '  -0.90%  '
    'D25' = '145.47'
    'E25' = '  +0.88%  '
    'E26' = '  +0.01%  '
    'E27' = '  -0.63%  '
    'E28' = '  +0.18%  '
    'D29' = '15.31'
    'E29' = '  +1.47%  '
    'E30' = '  -0.76%  '
    'E31' = '  +0.00%  '
    'E32' = '  +0.35%  '
    'D33' = '1.455.27'
    'E33' = '  +3.45%  '
    'E34' = '  +0.02%  '
    'D35' = '2.42'
    'E35' = '  -0.09%  '
    'E36' = '  +0.34%  '
    'E37' = '  -4.19%  '
    'E38' = '  -1.43%  '
    'D39' = '0.816'
    'E39' = '  -0.42%  '
    'D40' = '5.75'
    'E40' = '  -1.73%  '
    'E41' = '  +0.00%  '
    'E42' = '  +1.69%  '
    'D43' = '0.929'
    'E43' = '  -3.32%  '
    'D44' = '1.726.62'
    'E44' = '  +0.35%  '
    'D45' = '0.753'
    'E45' = '  -1.60%  '
    'D46' = '60.19'
    'E46' = '  -1.26%  '
    'D47' = '87.40'
    'E47' = '  +2.07%  '
    'B48' = 'BabyDogeCoin'
    'C48' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D48' = '0.0₆0102'
    'E48' = '  -2.06%  '
    'B49' = 'RenderToken'
    'C49' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D49' = '1.48'
    'E49' = '  -0.49%  '
    'E50' = '  +0.01%  '
    'D51' = '0.998'
    'E51' = '  +0.03%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
